$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows' timestamps / result values (RAD Phase 3 re-run values)
$ws.Range("B2").Value = "Thu Jan 25 17:46:39 EST 2024"
$ws.Range("B3").Value = "Thu Jan 25 17:46:52 EST 2024"
$ws.Range("B5").Value = "Thu Jan 25 17:47:05 EST 2024"

# Add two new test-data rows for Estate Tax
$ws.Range("D6").Value = "Existing Liability w/Notice Number"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").WrapText = $true
$ws.Range("D6").Borders.LineStyle = 1

$ws.Range("E6").Value = "Estate Tax"
$ws.Range("E6").WrapText = $true
$ws.Range("E6").Borders.LineStyle = 1

$ws.Range("D7").Value = "New Tax Return Amount Due"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").WrapText = $true
$ws.Range("D7").Borders.LineStyle = 1

$ws.Range("E7").Value = "Estate Tax"
$ws.Range("E7").WrapText = $true
$ws.Range("E7").Borders.LineStyle = 1

# Match the active selection left by the author
$ws.Range("E7").Select()
